# Update MSME indicator figures for Northern Mariana Islands with more
# precise (two decimal place) values, as per the 2015-04-01 autogenerated
# refresh of the source data.
#
# Row 11 - Enterprises density (per 1000 people): Micro/SMEs/MSMEs
# Row 12 - Employment (% of total): Micro/SMEs/MSMEs
# Row 14 - Enterprises (% of total): Micro/SMEs/MSMEs
#
# These cells hold their numeric-looking figures as plain text (shared
# strings), so we force the "Text" number format before writing the new
# value to keep Excel from re-interpreting the string as a number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "B11" = "12.91"
    "C11" = "9.53"
    "D11" = "22.44"
    "B12" = "10.48"
    "C12" = "33.39"
    "D12" = "43.87"
    "B14" = "28.76"
    "C14" = "37.94"
    "D14" = "89.32"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $originalStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    $cell.Style = $originalStyle
}
